$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.730734705924988
$ws.Range("B1").Value = 2.356780529022217
$ws.Range("C1").Value = 2.448511838912964
$ws.Range("D1").Value = 2.749891757965088
$ws.Range("E1").Value = 3.485612630844116
